$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("L2").Value = "[2.8910227675457527, 10.323382989350563]"
$ws.Range("M2").Value = 0.0005918973104912517
$ws.Range("N2").Value = 0.0005918973104912517
$ws.Range("R2").Value = 0.0002136552951506054
$ws.Range("T2").Value = "[7.109498362924011, 11.692784348550859]"
$ws.Range("U2").Value = 0.0000000000002220446049250313
$ws.Range("V2").Value = 0.0000000000002220446049250313

# Row 3 updates
$ws.Range("L3").Value = "[3.3316356372733757, 9.438232470904671]"
$ws.Range("M3").Value = 0.00005050883475488099
$ws.Range("N3").Value = 0.000101017669509762
$ws.Range("P3").Value = "[-1.3459476033397717, -0.4151053356094625]"
$ws.Range("Q3").Value = 0.000236457049741956
$ws.Range("R3").Value = 0.000236457049741956
$ws.Range("T3").Value = "[7.326373878006729, 10.4659529495893]"
$ws.Range("X3").Value = 1.541981981981996
$ws.Range("Y3").Value = 4.999759759759806
